$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.816.93'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('D3').Value = '2.498.59'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '322.99'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '108.76'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = "'0.550"
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').Value = '40.18'
$ws.Range('E10').Value = '  +4.94%  '
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '18.91'
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').Value = '7.21'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '2.890.00'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '2.496.58'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = "'0.850"
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').Value = '47.699.11'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').Value = '13.15'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('D20').Value = '6.63'
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0941'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = '2.77'
$ws.Range('E22').Value = '  +11.58%  '
$ws.Range('D23').Value = '70.74'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '247.86'
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  -1.21%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '25.91'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '9.97'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  -4.32%  '
$ws.Range('D30').Value = '35.19'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').Value = '19.85'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').Value = '5.35'
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('D35').Value = '0.0789'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').Value = '22.29'
$ws.Range('E41').Value = '  +5.52%  '
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').Value = '119.19'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').Value = '0.0297'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').Value = '1.999.76'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('D46').Value = '3.05'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('E47').Value = '  -3.46%  '
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('D49').Value = '9.02'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').Value = '5.11'
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('D51').Value = '56.92'
